# g_vs_hardcoded_COIN_6_checkers: refresh win-rate/avg-game-length table
# with results from 8 trained models (rows for checker counts 0-7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-5 (checker counts 0-3) ---
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 22

$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 3

$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 3

$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 2

# --- Append new rows 6-9 (checker counts 4-7) ---
# Column A holds text labels (e.g. "4"), so prefix with an apostrophe to
# force Excel to store them as text instead of numbers.
$ws.Range("A6").Value = "'4"
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 23

$ws.Range("A7").Value = "'5"
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 33

$ws.Range("A8").Value = "'6"
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 33

$ws.Range("A9").Value = "'7"
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 33

# Give the new label cells (A6:A9) the same formatting (bold, border,
# centered) as the existing label cells in column A, e.g. A2.
$ws.Range("A2").Copy()
$ws.Range("A6:A9").PasteSpecial(-4122)
